# Q4DemoSS2Metadata_v5_plainHeaders_new.xlsx
# "process renamed to protocol in many headers"
#
# Renames several shared-string header values from the old
# "*_process.process_core.*" / "process_*" naming convention to the new
# "*_protocol.protocol_core.*" / "protocol_*" naming convention on the
# "Enrichment protocol", "Library preparation protocol" and "Sequencing
# protocol" sheets (row 4 - the machine-readable header row), then restores
# the view/selection state that Excel recorded when the file was saved
# (per-sheet active cell, scroll position and the active tab).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Header (row 4) renames: process_core -> protocol_core, etc.
# ---------------------------------------------------------------------

# --- Enrichment protocol ---
$wsEnr = $wb.Worksheets.Item("Enrichment protocol")
$wsEnr.Range("B4").Value = "enrichment_protocol.protocol_core.protocol_name"
$wsEnr.Range("C4").Value = "enrichment_protocol.protocol_core.protocol_description"
$wsEnr.Range("D4").Value = "enrichment_protocol.protocol_core.start_time"
$wsEnr.Range("E4").Value = "enrichment_protocol.protocol_core.protocol_location"
$wsEnr.Range("F4").Value = "enrichment_protocol.protocol_core.operator_identity"

# --- Library preparation protocol ---
$wsLib = $wb.Worksheets.Item("Library preparation protocol")
$wsLib.Range("A4").Value = "library_preparation_protocol.protocol_core.protocol_id"
$wsLib.Range("B4").Value = "library_preparation_protocol.protocol_core.protocol_name"
$wsLib.Range("C4").Value = "library_preparation_protocol.protocol_core.protocol_description"
$wsLib.Range("D4").Value = "library_preparation_protocol.protocol_core.start_time"
$wsLib.Range("E4").Value = "library_preparation_protocol.protocol_core.protocol_location"
$wsLib.Range("F4").Value = "library_preparation_protocol.protocol_core.operator_identity"

# --- Sequencing protocol ---
$wsSeq = $wb.Worksheets.Item("Sequencing protocol")
$wsSeq.Range("A4").Value = "sequencing_protocol.protocol_core.protocol_id"
$wsSeq.Range("B4").Value = "sequencing_protocol.protocol_core.protocol_name"
$wsSeq.Range("C4").Value = "sequencing_protocol.protocol_core.protocol_description"
$wsSeq.Range("D4").Value = "sequencing_protocol.protocol_core.start_time"
$wsSeq.Range("E4").Value = "sequencing_protocol.protocol_core.protocol_location"
$wsSeq.Range("F4").Value = "sequencing_protocol.protocol_core.operator_identity"
$wsSeq.Range("L4").Value = "sequencing_protocol.protocol_type.text"

# ---------------------------------------------------------------------
# 2. View state: per-sheet active cell / scroll position, and which tab
#    ends up selected. Sheets are activated in the order Excel would have
#    left them so the final ActiveSheet / activeTab matches the saved file
#    (the last sheet activated below - "Sequence files" - is the one shown
#    when the workbook is reopened).
# ---------------------------------------------------------------------

$wsCellSusp = $wb.Worksheets.Item("Cell suspension")
$wsCellSusp.Activate()
$wsCellSusp.Range("AB4").Select()

$wsDiss = $wb.Worksheets.Item("Dissociation protocol")
$wsDiss.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 5
$wsDiss.Range("M4").Select()

$wsEnr.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
$wsEnr.Range("O4").Select()

$wsLib.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 34
$wsLib.Range("AO4").Select()

$wsSeq.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 7
$wsSeq.Range("M4").Select()

$wsFiles = $wb.Worksheets.Item("Sequence files")
$wsFiles.Activate()
$wsFiles.Range("M4").Select()
